# Update "想去人数" (F column) counts across the workbook's sheets to the
# freshly re-scraped values from the gh-pages data refresh at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets.Item(1))
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 12268
$ws1.Range("F3").Value  = 6899
$ws1.Range("F7").Value  = 263
$ws1.Range("F10").Value = 958
$ws1.Range("F11").Value = 117
$ws1.Range("F12").Value = 318
$ws1.Range("F13").Value = 966
$ws1.Range("F17").Value = 506
$ws1.Range("F18").Value = 213
$ws1.Range("F19").Value = 335
$ws1.Range("F21").Value = 249
$ws1.Range("F22").Value = 278
$ws1.Range("F23").Value = 67
$ws1.Range("F25").Value = 5109
$ws1.Range("F26").Value = 56
$ws1.Range("F27").Value = 1320
$ws1.Range("F29").Value = 785
$ws1.Range("F30").Value = 1270
$ws1.Range("F31").Value = 571

# Sheet "演出" (Worksheets.Item(2))
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value  = 3719
$ws2.Range("F10").Value = 1

# Sheet "本地生活" (Worksheets.Item(3))
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 1907

# Sheet "全部类型" (Worksheets.Item(4))
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 1907
$ws4.Range("F5").Value  = 12268
$ws4.Range("F6").Value  = 6899
$ws4.Range("F8").Value  = 3719
$ws4.Range("F12").Value = 263
$ws4.Range("F15").Value = 958
$ws4.Range("F16").Value = 117
$ws4.Range("F17").Value = 318
$ws4.Range("F18").Value = 966
$ws4.Range("F22").Value = 506
$ws4.Range("F23").Value = 213
$ws4.Range("F24").Value = 335
$ws4.Range("F26").Value = 249
$ws4.Range("F27").Value = 278
$ws4.Range("F33").Value = 5109
$ws4.Range("F34").Value = 56
$ws4.Range("F35").Value = 1320
$ws4.Range("F36").Value = 1
$ws4.Range("F40").Value = 785
$ws4.Range("F41").Value = 1270
$ws4.Range("F42").Value = 571
